$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 36 (shifts existing rows 36:139 down to 37:140)
$ws.Rows(36).Insert()

# Populate the newly inserted row 36 with the new record
$ws.Range("A36").Value = 11
$ws.Range("B36").Value = "Vega Monumental Concepción"
$ws.Range("C36").Value = "Bíobío"
$ws.Range("D36").Value = 44980
$ws.Range("E36").Value = 8
$ws.Range("F36").Value = 100112001
$ws.Range("G36").Value = "Berenjena"
$ws.Range("H36").Value = "Sin especificar"
$ws.Range("I36").Value = "Primera"
$ws.Range("J36").Value = 100
$ws.Range("K36").Value = 7000
$ws.Range("L36").Value = 7500
$ws.Range("M36").Value = 7250
$ws.Range("N36").Value = "`$/caja 60 unidades"
$ws.Range("O36").Value = "Región de Arica y Parinacota"
$ws.Range("P36").Value = 121
$ws.Range("Q36").Value = 60
$ws.Range("R36").Value = "Hortaliza"
